# Update "paises" (countries) COVID dashboard data.
# 1) Refresh the "last updated" timestamp in the title cell.
# 2) Apply the updated case/death/recovery figures for the countries whose
#    numbers changed since the last refresh.
# 3) Re-sort the country table (rows 4:219) by "Casos totales" (column B)
#    descending, as the live dashboard does on every refresh - this is what
#    naturally reshuffles countries with tied/close totals (e.g. San Marino
#    and Republica de Africa Central; Groenlandia and Islas Turcas y Caicos).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Update timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 26 de Mayo de 2020 a las 19:05"

# --- 2) Update country figures -------------------------------------------
# Estados Unidos (row 4)
$ws.Cells.Item(4, 2).Value = 1713463
$ws.Cells.Item(4, 3).Value = 7237
$ws.Cells.Item(4, 4).Value = 468778
$ws.Cells.Item(4, 5).Value = 1144644
$ws.Cells.Item(4, 7).Value = 236
$ws.Cells.Item(4, 8).Value = 100041

# Alemania (row 11)
$ws.Cells.Item(11, 2).Value = 181094
$ws.Cells.Item(11, 3).Value = 305
$ws.Cells.Item(11, 5).Value = 10645
$ws.Cells.Item(11, 7).Value = 21
$ws.Cells.Item(11, 8).Value = 8449

# India (row 13)
$ws.Cells.Item(13, 2).Value = 150739
$ws.Cells.Item(13, 3).Value = 5789
$ws.Cells.Item(13, 4).Value = 64272
$ws.Cells.Item(13, 5).Value = 82118
$ws.Cells.Item(13, 7).Value = 177
$ws.Cells.Item(13, 8).Value = 4349

# Canada (row 16)
$ws.Cells.Item(16, 2).Value = 85998
$ws.Cells.Item(16, 3).Value = 287
$ws.Cells.Item(16, 5).Value = 34534

# Rumania (row 40)
$ws.Cells.Item(40, 5).Value = 5339
$ws.Cells.Item(40, 7).Value = 11
$ws.Cells.Item(40, 8).Value = 1216

# Israel (row 42)
$ws.Cells.Item(42, 2).Value = 16757
$ws.Cells.Item(42, 3).Value = 23
$ws.Cells.Item(42, 4).Value = 14457
$ws.Cells.Item(42, 5).Value = 2019

# Argelia (row 56)
$ws.Cells.Item(56, 2).Value = 8697
$ws.Cells.Item(56, 3).Value = 194
$ws.Cells.Item(56, 4).Value = 4918
$ws.Cells.Item(56, 5).Value = 3162
$ws.Cells.Item(56, 7).Value = 8
$ws.Cells.Item(56, 8).Value = 617

# Paraguay (row 118)
$ws.Cells.Item(118, 2).Value = 877
$ws.Cells.Item(118, 3).Value = 12
$ws.Cells.Item(118, 4).Value = 382
$ws.Cells.Item(118, 5).Value = 484

# Republica del Chad (row 129)
$ws.Cells.Item(129, 2).Value = 700
$ws.Cells.Item(129, 3).Value = 13
$ws.Cells.Item(129, 4).Value = 303
$ws.Cells.Item(129, 5).Value = 335
$ws.Cells.Item(129, 7).Value = 1
$ws.Cells.Item(129, 8).Value = 62

# Republica de Africa Central (currently at row 131, before sorting)
# San Marino (currently at row 130) is unchanged and simply gets outranked
# once Africa Central's total overtakes it, so only row 131 needs updating.
$ws.Cells.Item(131, 2).Value = 671
$ws.Cells.Item(131, 3).Value = 19
$ws.Cells.Item(131, 5).Value = 648

# Suazilandia (row 152)
$ws.Cells.Item(152, 2).Value = 261
$ws.Cells.Item(152, 3).Value = 5
$ws.Cells.Item(152, 4).Value = 164
$ws.Cells.Item(152, 5).Value = 95

# Islas Turcas y Caicos (row 207) and Groenlandia (row 208) have identical
# totals (12) both before and after, so a totals-based sort leaves them tied
# and keeps their original relative order. The refreshed source data lists
# them in the opposite order, so swap the two whole rows explicitly (their
# underlying per-country numbers are otherwise unchanged).
$turcas = $ws.Range("A207:H207").Value()
$groenlandia = $ws.Range("A208:H208").Value()
$ws.Range("A207:H207").Value = $groenlandia
$ws.Range("A208:H208").Value = $turcas

# Sahara Occidental (row 212)
$ws.Cells.Item(212, 5).Value = 2
$ws.Cells.Item(212, 7).Value = 1
$ws.Cells.Item(212, 8).Value = 1

# --- 3) Re-sort the table by Casos totales (column B), descending --------
$dataRange = $ws.Range("A4:H219")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("B4:B219"), 0, 2, 0, 0)
$ws.Sort.SetRange($dataRange)
$ws.Sort.Header = 2
$ws.Sort.Apply()
